# Updated cryptos list on Thu Jan 25 19:26:07 UTC 2024 with GitHub Actions
#
# Refresh the Price (D) and Volume(1h) (E) columns for the cryptos table,
# plus reorder a couple of coin pairs that swapped rank (Toncoin/Cosmos,
# VeChain/FraxShare, ApeXProtocol/EnergySwap) so B/C/D/E all move together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells whose new value is a "pure" number-looking string (e.g.
# "9.23") that Excel's normal text coercion would otherwise snap into a
# numeric cell. Force these ranges to Text format before writing so they
# stay inline/shared strings exactly like the surrounding cells, then
# restore the default "Normal" style so no stray number-format survives
# on the cell (matches the original workbook, which has no explicit style
# on any of these cells).
$textForceAddrs = @(
    "D5","D6","D9","D10","D11","D12","D14","D16","D18","D21","D22","D23",
    "D24","D26","D27","D28","D29","D30","D31","D32","D34","D35","D36",
    "D39","D40","D41","D43","D44","D45","D46","D47","D51"
)

foreach ($addr in $textForceAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Price (D) / Volume(1h) (E) refresh, row by row ---

# Row 2: Bitcoin
$ws.Range("D2").Value = "39.871.74"
$ws.Range("E2").Value = "  +0.17%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.212.88"
$ws.Range("E3").Value = "  -0.18%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5: BNB
$ws.Range("D5").Value = "291.81"
$ws.Range("E5").Value = "  -0.23%  "

# Row 6: Solana
$ws.Range("D6").Value = "87.09"
$ws.Range("E6").Value = "  +1.32%  "

# Row 7: XRP
$ws.Range("E7").Value = "  -0.11%  "

# Row 8: USDC
$ws.Range("E8").Value = "  -0.07%  "

# Row 9: Cardano
$ws.Range("D9").Value = "0.468"
$ws.Range("E9").Value = "  -0.78%  "

# Row 10: Avalanche
$ws.Range("D10").Value = "30.35"
$ws.Range("E10").Value = "  -0.95%  "

# Row 11: Dogecoin
$ws.Range("D11").Value = "0.0780"
$ws.Range("E11").Value = "  -0.65%  "

# Row 12: OKB
$ws.Range("D12").Value = "49.98"
$ws.Range("E12").Value = "  +5.45%  "

# Row 13: TRON
$ws.Range("E13").Value = "  +2.56%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "6.45"
$ws.Range("E14").Value = "  +1.69%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.551.52"
$ws.Range("E15").Value = "  -0.38%  "

# Row 16: Chainlink
$ws.Range("D16").Value = "13.74"
$ws.Range("E16").Value = "  -2.02%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "2.216.93"
$ws.Range("E17").Value = "  +0.15%  "

# Row 18: Polygon
$ws.Range("D18").Value = "0.731"
$ws.Range("E18").Value = "  +0.04%  "

# Row 19: WrappedBTC
$ws.Range("D19").Value = "39.790.85"
$ws.Range("E19").Value = "  +0.03%  "

# Row 20: ShibaInu
$ws.Range("D20").Value = "0.0₃0885"
$ws.Range("E20").Value = "  +0.41%  "

# Row 21: InternetComputer(DFINITY)
$ws.Range("D21").Value = "11.16"
$ws.Range("E21").Value = "  +0.16%  "

# Row 22: Uniswap
$ws.Range("D22").Value = "5.74"
$ws.Range("E22").Value = "  -0.97%  "

# Row 23: Litecoin
$ws.Range("D23").Value = "65.54"
$ws.Range("E23").Value = "  +0.13%  "

# Row 24: BitcoinCash
$ws.Range("D24").Value = "237.27"
$ws.Range("E24").Value = "  +0.62%  "

# Row 25: Dai
$ws.Range("E25").Value = "  +0.14%  "

# Row 26: PancakeSwap
$ws.Range("D26").Value = "2.46"
$ws.Range("E26").Value = "  -0.26%  "

# Row 27: ImmutableX
$ws.Range("D27").Value = "1.83"
$ws.Range("E27").Value = "  -0.18%  "

# Row 28: EthereumClassic
$ws.Range("D28").Value = "23.17"
$ws.Range("E28").Value = "  +1.98%  "

# Row 29 & 30: Toncoin and Cosmos swap places (row 29 becomes Cosmos,
# row 30 becomes Toncoin), each carrying its own refreshed price/volume.
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "9.23"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.05"
$ws.Range("E30").Value = "  -6.95%  "

# Row 31: Monero
$ws.Range("D31").Value = "157.23"
$ws.Range("E31").Value = "  +3.67%  "

# Row 32: InjectiveProtocol
$ws.Range("D32").Value = "31.92"
$ws.Range("E32").Value = "  -2.72%  "

# Row 33: FirstDigitalUSD
$ws.Range("E33").Value = "  -0.04%  "

# Row 34: Filecoin
$ws.Range("D34").Value = "4.96"
$ws.Range("E34").Value = "  +0.58%  "

# Row 35: Hedera
$ws.Range("D35").Value = "0.0710"
$ws.Range("E35").Value = "  -1.15%  "

# Row 36: LidoDAOToken
$ws.Range("D36").Value = "2.93"
$ws.Range("E36").Value = "  +4.93%  "

# Row 37: WEMIXToken
$ws.Range("E37").Value = "  -1.49%  "

# Row 38: Stellar
$ws.Range("E38").Value = "  -0.48%  "

# Row 39: Kaspa
$ws.Range("D39").Value = "0.0982"
$ws.Range("E39").Value = "  -0.94%  "

# Row 40: ARBITRUM
$ws.Range("D40").Value = "1.72"
$ws.Range("E40").Value = "  +1.24%  "

# Row 41: Celestia
$ws.Range("D41").Value = "15.25"
$ws.Range("E41").Value = "  -3.94%  "

# Row 42: Maker
$ws.Range("D42").Value = "2.113.58"
$ws.Range("E42").Value = "  +2.46%  "

# Row 43: RenderToken
$ws.Range("D43").Value = "3.73"
$ws.Range("E43").Value = "  -1.41%  "

# Row 44 & 45: VeChain and FraxShare swap places (row 44 becomes
# FraxShare, row 45 becomes VeChain).
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "10.03"
$ws.Range("E44").Value = "  +1.22%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0269"
$ws.Range("E45").Value = "  +0.55%  "

# Row 46 & 47: ApeXProtocol and EnergySwap swap places (row 46 becomes
# EnergySwap, row 47 becomes ApeXProtocol).
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "17.85"
$ws.Range("E46").Value = "  +0.38%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "2.09"
$ws.Range("E47").Value = "  -0.28%  "

# Row 48: NEARProtocol
$ws.Range("E48").Value = "  +3.73%  "

# Row 49: RocketPoolETH
$ws.Range("D49").Value = "2.422.75"

# Row 50: Stacks
$ws.Range("E50").Value = "  +3.08%  "

# Row 51: Aave
$ws.Range("D51").Value = "88.57"
$ws.Range("E51").Value = "  -0.37%  "

# Restore the default cell style on every cell we force-formatted to Text,
# so the saved XML carries no explicit style index (matching the rest of
# the untouched data cells in this sheet).
foreach ($addr in $textForceAddrs) {
    $ws.Range($addr).Style = "Normal"
}
